# Apply "repull data, push all data, mean calculation" edits:
# Update the dSF column (F) values for the rows whose recalculated
# values changed after repulling the source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "F4"  = -3
    "F5"  = -1
    "F6"  = 4
    "F7"  = 2
    "F9"  = 2
    "F12" = -5
    "F13" = -8
    "F14" = -4
    "F16" = -5
    "F20" = 2
    "F22" = -6
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
